$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.824.30'
$ws.Range("E2").Value = '  -0.23%  '

$ws.Range("D3").Value = '2.220.21'
$ws.Range("E3").Value = '  -0.12%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '253.01'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +5.20%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.631'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.98%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '70.44'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +4.34%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.603'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +10.39%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.62'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +12.75%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0969'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.83%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '58.44'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.66%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.34'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +10.79%  '

$ws.Range("E14").Value = '  +0.14%  '

$ws.Range("D15").Value = '2.547.17'
$ws.Range("E15").Value = '  -0.54%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.05'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.22%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.897'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +5.95%  '

$ws.Range("D18").Value = '2.212.17'
$ws.Range("E18").Value = '  -0.70%  '

$ws.Range("D19").Value = '41.736.03'
$ws.Range("E19").Value = '  -0.21%  '

$ws.Range("D20").Value = '0.0₃0967'
$ws.Range("E20").Value = '  +1.77%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.30'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +3.13%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.72'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.42%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.45'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.90%  '

$ws.Range("E24").Value = '  +1.45%  '

$ws.Range("E25").Value = '  +13.38%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.09'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +23.38%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.06%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.53'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +4.41%  '

$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.20'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.34%  '

$ws.Range("B30").Value = 'Monero'
$ws.Range("C30").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '171.83'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.60%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.86'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +2.94%  '

$ws.Range("E32").Value = '  +2.79%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.65'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +9.86%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.125'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.13%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0748'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +5.78%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.70'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.38%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '25.84'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +12.13%  '

$ws.Range("E38").Value = '  +5.98%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0306'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +10.20%  '

$ws.Range("E40").Value = '  +0.41%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.94'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +2.59%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '12.38'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +28.86%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '65.22'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.04%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.206'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +9.00%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.87'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.72%  '

$ws.Range("E46").Value = '  -1.39%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.103'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.40%  '

$ws.Range("E48").Value = '  +0.17%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.57'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.79%  '

$ws.Range("B50").Value = 'ARBITRUM'
$ws.Range("C50").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.17'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +7.07%  '

$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.43'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +7.02%  '
